$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.3327402135231317
$wsSummary.Range("C2").Value = 0.06516290726817042
$wsSummary.Range("D2").Value = 0.9285714285714286
$wsSummary.Range("E2").Value = 0.1217798594847775
$wsSummary.Range("F2").Value = 0.2544031311154599
$wsSummary.Range("G2").Value = 0.6151046405823476
$wsSummary.Range("H2").Value = 0.8103263777421081
$wsSummary.Range("I2").Value = 26
$wsSummary.Range("J2").Value = 373
$wsSummary.Range("K2").Value = 161
$wsSummary.Range("L2").Value = 2

# --- Classification Report sheet ---
$wsClassification = $wb.Worksheets.Item("Classification Report")
$wsClassification.Range("B2").Value = 0.9877300613496932
$wsClassification.Range("C2").Value = 0.301498127340824
$wsClassification.Range("D2").Value = 0.4619799139167862

$wsClassification.Range("B3").Value = 0.06516290726817042
$wsClassification.Range("C3").Value = 0.9285714285714286
$wsClassification.Range("D3").Value = 0.1217798594847775

$wsClassification.Range("B4").Value = 0.3327402135231317
$wsClassification.Range("C4").Value = 0.3327402135231317
$wsClassification.Range("D4").Value = 0.3327402135231317
$wsClassification.Range("E4").Value = 0.3327402135231317

$wsClassification.Range("B5").Value = 0.5264464843089318
$wsClassification.Range("C5").Value = 0.6150347779561263
$wsClassification.Range("D5").Value = 0.2918798867007819

$wsClassification.Range("B6").Value = 0.9417658615022153
$wsClassification.Range("C6").Value = 0.3327402135231317
$wsClassification.Range("D6").Value = 0.4450304450127003

# --- Confusion Matrix sheet ---
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")
$wsConfusion.Range("B2").Value = 161
$wsConfusion.Range("C2").Value = 373
$wsConfusion.Range("B3").Value = 2
$wsConfusion.Range("C3").Value = 26
